$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 636 (shifts existing rows 636..677 down to 637..678)
$ws.Rows.Item(636).Insert()

# Populate the newly inserted row with the new data point
# Force text format on column A so the date-like string isn't coerced to a date serial
$ws.Range("A636").NumberFormat = "@"
$ws.Range("A636").Value = "2026/01/14"
$ws.Range("B636").Value = "水"
$ws.Range("C636").Value = 14
$ws.Range("D636").Value = 201
